$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price column (D) so numeric-looking strings
# (e.g. "15.30") keep their exact text representation instead of being
# auto-converted to numbers (which would drop trailing zeros, etc.).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '37.269.81'
$ws.Range('E2').Value = '  +2.30%  '
$ws.Range('D3').Value = '2.094.77'
$ws.Range('E3').Value = '  +4.40%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '250.42'
$ws.Range('E5').Value = '  +2.34%  '
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '53.66'
$ws.Range('E8').Value = '  +20.21%  '
$ws.Range('D9').Value = '61.79'
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('D10').Value = '0.374'
$ws.Range('E10').Value = '  +1.75%  '
$ws.Range('D11').Value = '0.0741'
$ws.Range('E11').Value = '  +4.16%  '
$ws.Range('E12').Value = '  +7.76%  '
$ws.Range('D13').Value = '15.30'
$ws.Range('E13').Value = '  +5.12%  '
$ws.Range('D14').Value = '2.400.37'
$ws.Range('E14').Value = '  +4.35%  '
$ws.Range('D15').Value = '0.837'
$ws.Range('E15').Value = '  +3.35%  '
$ws.Range('D16').Value = '2.097.50'
$ws.Range('E16').Value = '  +4.42%  '
$ws.Range('D17').Value = '5.16'
$ws.Range('E17').Value = '  +5.75%  '
$ws.Range('D18').Value = '37.186.47'
$ws.Range('E18').Value = '  +2.22%  '
$ws.Range('D19').Value = '72.77'
$ws.Range('E19').Value = '  +2.10%  '
$ws.Range('D20').Value = '14.63'
$ws.Range('E20').Value = '  +14.21%  '
$ws.Range('D21').Value = '0.0₃0841'
$ws.Range('E21').Value = '  +3.38%  '
$ws.Range('D22').Value = '240.98'
$ws.Range('E22').Value = '  +1.67%  '
$ws.Range('D23').Value = '5.20'
$ws.Range('E23').Value = '  +6.70%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('D26').Value = '171.92'
$ws.Range('E26').Value = '  +3.86%  '
$ws.Range('D27').Value = '9.29'
$ws.Range('E27').Value = '  +8.11%  '
$ws.Range('D28').Value = '20.65'
$ws.Range('E28').Value = '  +5.52%  '
$ws.Range('D29').Value = '1.99'
$ws.Range('E29').Value = '  +3.63%  '
$ws.Range('E30').Value = '  +1.56%  '
$ws.Range('E31').Value = '  +26.52%  '
$ws.Range('D32').Value = '22.15'
$ws.Range('E32').Value = '  -0.69%  '
$ws.Range('E33').Value = '  +3.18%  '
$ws.Range('D34').Value = '0.0615'
$ws.Range('E34').Value = '  +5.40%  '
$ws.Range('D35').Value = '0.0905'
$ws.Range('E35').Value = '  +12.03%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').Value = '2.28'
$ws.Range('E37').Value = '  +6.78%  '
$ws.Range('D38').Value = '4.11'
$ws.Range('E38').Value = '  +3.39%  '
$ws.Range('D39').Value = '1.84'
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('D40').Value = '1.33'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').Value = '18.37'
$ws.Range('E41').Value = '  +14.93%  '
$ws.Range('E42').Value = '  +4.50%  '
$ws.Range('E43').Value = '  +5.70%  '
$ws.Range('D44').Value = '98.68'
$ws.Range('E44').Value = '  +3.20%  '
$ws.Range('D45').Value = '0.0923'
$ws.Range('E45').Value = '  +13.10%  '
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('D47').Value = '4.11'
$ws.Range('E47').Value = '  +99.14%  '
$ws.Range('D48').Value = '1.317.21'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('D49').Value = '2.95'
$ws.Range('E49').Value = '  +7.08%  '
$ws.Range('D50').Value = '7.02'
$ws.Range('E50').Value = '  +14.27%  '
$ws.Range('D51').Value = '2.284.19'
$ws.Range('E51').Value = '  +4.08%  '
